$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "2025-04-28 21:38:00"
$ws.Range("B21").Value = 42
